# "minor fix for text after grouping"
# Adds three "Text after" rows following the existing @group/@endgroup
# template block (rows 1-5) on Sheet1, and moves the selection down to
# the first empty row beneath the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = "Text after"
$ws.Range("A7").Value = "Text after"
$ws.Range("A8").Value = "Text after"

$ws.Range("A9").Select()
